# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-10-05 (serial 45204) to 2023-10-06 (serial 45205).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162, find last used row in column A
if ($lastRow -lt 2) { $lastRow = 399 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45205
